# Chicker battery_offset.xlsx - replace simulated extrapolated data with
# real measured data points (9 rows instead of 11), drop the formula-driven
# extrapolation in columns A/C in favor of literal measured values, and add
# the linear-fit trendline equation text below the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New measured "Actual Voltage" (A) / "Simulated Voltage" (C) pairs ----
$ws.Range("A2").Value2 = 3.1
$ws.Range("A3").Value2 = 4.45
$ws.Range("A4").Value2 = 6.35
$ws.Range("A5").Value2 = 8.06
$ws.Range("A6").Value2 = 10.64
$ws.Range("A7").Value2 = 13
$ws.Range("A8").Value2 = 14.82
$ws.Range("A9").Value2 = 15.98

$ws.Range("C2").Value2 = 2.8
$ws.Range("C3").Value2 = 4.3
$ws.Range("C4").Value2 = 6.3
$ws.Range("C5").Value2 = 8.2
$ws.Range("C6").Value2 = 11
$ws.Range("C7").Value2 = 13.6
$ws.Range("C8").Value2 = 15.6
$ws.Range("C9").Value2 = 16.8

# Column B keeps the calibration formula A/1.00614, now only spanning the
# 8 data rows (2-9) as one shared formula.
$ws.Range("B2:B9").Formula = "=A2/1.00614"

# Rows 10 and 11 (old extrapolated points) are gone entirely, and the old
# row 12 data (A12/B12/C12) is gone too - only the trendline equation
# label remains, moved into D12.
$ws.Range("A10:C11").Clear()
$ws.Range("A12:C12").Clear()
$ws.Range("D12").Value2 = "y = 1.0893x - 0.5781"

# Match the saved selection state.
$ws.Range("D12").Select() | Out-Null
